$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 10891  # ALC!H43
$ws.Cells.Item(43, 10).Value = 983.5  # ALC!J43
$ws.Cells.Item(43, 12).Value = 983.5  # ALC!L43
$ws.Cells.Item(43, 14).Value = -1121.5  # ALC!N43

$ws.Cells.Item(62, 8).Value = 10354.896  # ALC!H62
$ws.Cells.Item(62, 9).Value = 8330.23  # ALC!I62
$ws.Cells.Item(62, 11).Value = 8330.23  # ALC!K62
$ws.Cells.Item(62, 13).Value = -7706.23  # ALC!M62

$ws.Cells.Item(65, 8).Value = 10354.896  # ALC!H65
$ws.Cells.Item(65, 9).Value = 8330.23  # ALC!I65
$ws.Cells.Item(65, 11).Value = 41651.14999999999  # ALC!K65
$ws.Cells.Item(65, 13).Value = -38531.14999999999  # ALC!M65

$ws.Cells.Item(94, 8).Value = 29391.234  # ALC!H94
$ws.Cells.Item(94, 9).Value = 30040.625  # ALC!I94
$ws.Cells.Item(94, 11).Value = 30040.625  # ALC!K94
$ws.Cells.Item(94, 13).Value = -29589.625  # ALC!M94

$ws.Cells.Item(96, 8).Value = 913.7143  # ALC!H96
$ws.Cells.Item(96, 9).Value = 1045.9333  # ALC!I96
$ws.Cells.Item(96, 11).Value = 3137.7999  # ALC!K96
$ws.Cells.Item(96, 13).Value = -1764.7999  # ALC!M96

$ws.Cells.Item(111, 8).Value = 1449.8334  # ALC!H111
$ws.Cells.Item(111, 9).Value = 1140  # ALC!I111
$ws.Cells.Item(111, 11).Value = 3420  # ALC!K111
$ws.Cells.Item(111, 13).Value = -353  # ALC!M111

$ws.Cells.Item(132, 8).Value = 4220341.5  # ALC!H132
$ws.Cells.Item(132, 9).Value = 4831740  # ALC!I132
$ws.Cells.Item(132, 11).Value = 14495220  # ALC!K132
$ws.Cells.Item(132, 13).Value = -14492690  # ALC!M132

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3316.0166  # ARM!H32
$ws.Cells.Item(32, 9).Value = 2790.7544  # ARM!I32
$ws.Cells.Item(32, 11).Value = 2790.7544  # ARM!K32
$ws.Cells.Item(32, 13).Value = -2503.7544  # ARM!M32

$ws.Cells.Item(45, 8).Value = 6196.607  # ARM!H45
$ws.Cells.Item(45, 9).Value = 6241.1113  # ARM!I45
$ws.Cells.Item(45, 10).Value = 4995  # ARM!J45
$ws.Cells.Item(45, 11).Value = 6241.1113  # ARM!K45
$ws.Cells.Item(45, 12).Value = 4995  # ARM!L45
$ws.Cells.Item(45, 13).Value = -5864.1113  # ARM!M45
$ws.Cells.Item(45, 14).Value = -5749  # ARM!N45

$ws.Cells.Item(61, 8).Value = 3811.9614  # ARM!H61
$ws.Cells.Item(61, 9).Value = 3084.7368  # ARM!I61
$ws.Cells.Item(61, 11).Value = 3084.7368  # ARM!K61
$ws.Cells.Item(61, 13).Value = -2872.7368  # ARM!M61

$ws.Cells.Item(74, 8).Value = 5148.6304  # ARM!H74
$ws.Cells.Item(74, 9).Value = 1183.0256  # ARM!I74
$ws.Cells.Item(74, 11).Value = 1183.0256  # ARM!K74
$ws.Cells.Item(74, 13).Value = -309.0255999999999  # ARM!M74

$ws.Cells.Item(77, 8).Value = 5148.6304  # ARM!H77
$ws.Cells.Item(77, 9).Value = 1183.0256  # ARM!I77
$ws.Cells.Item(77, 11).Value = 5915.128  # ARM!K77
$ws.Cells.Item(77, 13).Value = -1547.128  # ARM!M77

$ws.Cells.Item(97, 8).Value = 829.45  # ARM!H97
$ws.Cells.Item(97, 10).Value = 105  # ARM!J97
$ws.Cells.Item(97, 12).Value = 105  # ARM!L97
$ws.Cells.Item(97, 14).Value = -1097  # ARM!N97

$ws.Cells.Item(102, 8).Value = 3858.4783  # ARM!H102
$ws.Cells.Item(102, 9).Value = 3086.8667  # ARM!I102
$ws.Cells.Item(102, 11).Value = 3086.8667  # ARM!K102
$ws.Cells.Item(102, 13).Value = -1464.8667  # ARM!M102

$ws.Cells.Item(110, 8).Value = 6196.185  # ARM!H110
$ws.Cells.Item(110, 9).Value = 7072.1113  # ARM!I110
$ws.Cells.Item(110, 11).Value = 7072.1113  # ARM!K110
$ws.Cells.Item(110, 13).Value = -5027.1113  # ARM!M110

$ws.Cells.Item(132, 8).Value = 2569.8823  # ARM!H132
$ws.Cells.Item(132, 9).Value = 1993.75  # ARM!I132
$ws.Cells.Item(132, 11).Value = 5981.25  # ARM!K132
$ws.Cells.Item(132, 13).Value = -3451.25  # ARM!M132

$ws.Cells.Item(136, 8).Value = 3811.9614  # ARM!H136
$ws.Cells.Item(136, 9).Value = 3084.7368  # ARM!I136
$ws.Cells.Item(136, 11).Value = 9254.2104  # ARM!K136
$ws.Cells.Item(136, 13).Value = -6704.2104  # ARM!M136

$ws.Cells.Item(139, 8).Value = 50357.5  # ARM!H139
$ws.Cells.Item(139, 10).Value = 50357.5  # ARM!J139
$ws.Cells.Item(139, 12).Value = 50357.5  # ARM!L139
$ws.Cells.Item(139, 14).Value = -60637.5  # ARM!N139

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(81, 8).Value = 37390  # BSM!H81
$ws.Cells.Item(81, 10).Value = 37390  # BSM!J81
$ws.Cells.Item(81, 12).Value = 37390  # BSM!L81
$ws.Cells.Item(81, 14).Value = -39512  # BSM!N81

$ws.Cells.Item(84, 8).Value = 37390  # BSM!H84
$ws.Cells.Item(84, 10).Value = 37390  # BSM!J84
$ws.Cells.Item(84, 12).Value = 112170  # BSM!L84
$ws.Cells.Item(84, 14).Value = -122778  # BSM!N84

$ws.Cells.Item(86, 8).Value = 3604.923  # BSM!H86
$ws.Cells.Item(86, 9).Value = 910.6667  # BSM!I86
$ws.Cells.Item(86, 10).Value = 5914.2856  # BSM!J86
$ws.Cells.Item(86, 11).Value = 910.6667  # BSM!K86
$ws.Cells.Item(86, 12).Value = 5914.2856  # BSM!L86
$ws.Cells.Item(86, 13).Value = 212.3333  # BSM!M86
$ws.Cells.Item(86, 14).Value = -8160.2856  # BSM!N86

$ws.Cells.Item(89, 8).Value = 3604.923  # BSM!H89
$ws.Cells.Item(89, 9).Value = 910.6667  # BSM!I89
$ws.Cells.Item(89, 10).Value = 5914.2856  # BSM!J89
$ws.Cells.Item(89, 11).Value = 4553.3335  # BSM!K89
$ws.Cells.Item(89, 12).Value = 29571.428  # BSM!L89
$ws.Cells.Item(89, 13).Value = 1062.6665  # BSM!M89
$ws.Cells.Item(89, 14).Value = -40803.428  # BSM!N89

$ws.Cells.Item(134, 8).Value = 1910.7435  # BSM!H134
$ws.Cells.Item(134, 9).Value = 1908.3948  # BSM!I134
$ws.Cells.Item(134, 11).Value = 5725.1844  # BSM!K134
$ws.Cells.Item(134, 13).Value = -3190.1844  # BSM!M134

$ws.Cells.Item(135, 8).Value = 55000  # BSM!H135
$ws.Cells.Item(135, 10).Value = 55000  # BSM!J135
$ws.Cells.Item(135, 12).Value = 55000  # BSM!L135
$ws.Cells.Item(135, 14).Value = -65140  # BSM!N135

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 72329.13  # CRP!H31
$ws.Cells.Item(31, 9).Value = 113252  # CRP!I31
$ws.Cells.Item(31, 11).Value = 113252  # CRP!K31
$ws.Cells.Item(31, 13).Value = -112957  # CRP!M31

$ws.Cells.Item(34, 8).Value = 72329.13  # CRP!H34
$ws.Cells.Item(34, 9).Value = 113252  # CRP!I34
$ws.Cells.Item(34, 11).Value = 113252  # CRP!K34
$ws.Cells.Item(34, 13).Value = -113050  # CRP!M34

$ws.Cells.Item(75, 8).Value = 15000  # CRP!H75
$ws.Cells.Item(75, 10).Value = 15000  # CRP!J75
$ws.Cells.Item(75, 12).Value = 15000  # CRP!L75
$ws.Cells.Item(75, 14).Value = -16996  # CRP!N75

$ws.Cells.Item(78, 8).Value = 15000  # CRP!H78
$ws.Cells.Item(78, 10).Value = 15000  # CRP!J78
$ws.Cells.Item(78, 12).Value = 45000  # CRP!L78
$ws.Cells.Item(78, 14).Value = -54984  # CRP!N78

$ws.Cells.Item(134, 8).Value = 9157.026  # CRP!H134
$ws.Cells.Item(134, 9).Value = 5843.241  # CRP!I134
$ws.Cells.Item(134, 10).Value = 19834.777  # CRP!J134
$ws.Cells.Item(134, 11).Value = 17529.723  # CRP!K134
$ws.Cells.Item(134, 12).Value = 59504.33099999999  # CRP!L134
$ws.Cells.Item(134, 13).Value = -14994.723  # CRP!M134
$ws.Cells.Item(134, 14).Value = -64574.33099999999  # CRP!N134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 9).Value = 48.565216  # CUL!I2
$ws.Cells.Item(2, 10).Value = 8.333333  # CUL!J2
$ws.Cells.Item(2, 11).Value = 291.391296  # CUL!K2
$ws.Cells.Item(2, 12).Value = 49.999998  # CUL!L2
$ws.Cells.Item(2, 13).Value = -178.391296  # CUL!M2
$ws.Cells.Item(2, 14).Value = -275.999998  # CUL!N2

$ws.Cells.Item(38, 8).Value = 120.125  # CUL!H38
$ws.Cells.Item(38, 9).Value = 69.5  # CUL!I38
$ws.Cells.Item(38, 11).Value = 208.5  # CUL!K38
$ws.Cells.Item(38, 13).Value = 138.5  # CUL!M38

$ws.Cells.Item(75, 8).Value = 9499.333000000001  # CUL!H75
$ws.Cells.Item(75, 10).Value = 10499  # CUL!J75
$ws.Cells.Item(75, 12).Value = 31497  # CUL!L75
$ws.Cells.Item(75, 14).Value = -33493  # CUL!N75

$ws.Cells.Item(78, 8).Value = 9499.333000000001  # CUL!H78
$ws.Cells.Item(78, 10).Value = 10499  # CUL!J78
$ws.Cells.Item(78, 12).Value = 94491  # CUL!L78
$ws.Cells.Item(78, 14).Value = -104475  # CUL!N78

$ws.Cells.Item(94, 8).Value = 5007.5835  # CUL!H94
$ws.Cells.Item(94, 9).Value = 2749.5  # CUL!I94
$ws.Cells.Item(94, 11).Value = 8248.5  # CUL!K94
$ws.Cells.Item(94, 13).Value = -7572.5  # CUL!M94

$ws.Cells.Item(131, 8).Value = 50395.094  # CUL!H131
$ws.Cells.Item(131, 10).Value = 2914.85  # CUL!J131
$ws.Cells.Item(131, 12).Value = 8744.549999999999  # CUL!L131
$ws.Cells.Item(131, 14).Value = -18824.55  # CUL!N131

$ws.Cells.Item(132, 8).Value = 1222.2222  # CUL!H132
$ws.Cells.Item(132, 9).Value = 1000  # CUL!I132
$ws.Cells.Item(132, 10).Value = 1666.6666  # CUL!J132
$ws.Cells.Item(132, 11).Value = 9000  # CUL!K132
$ws.Cells.Item(132, 12).Value = 14999.9994  # CUL!L132
$ws.Cells.Item(132, 13).Value = -6470  # CUL!M132
$ws.Cells.Item(132, 14).Value = -20059.9994  # CUL!N132

$ws.Cells.Item(137, 8).Value = 2624.762  # CUL!H137
$ws.Cells.Item(137, 9).Value = 1431.3636  # CUL!I137
$ws.Cells.Item(137, 10).Value = 3937.5  # CUL!J137
$ws.Cells.Item(137, 11).Value = 4294.0908  # CUL!K137
$ws.Cells.Item(137, 12).Value = 11812.5  # CUL!L137
$ws.Cells.Item(137, 13).Value = 805.9092000000001  # CUL!M137
$ws.Cells.Item(137, 14).Value = -22012.5  # CUL!N137

$ws.Cells.Item(141, 8).Value = 257233.25  # CUL!H141
$ws.Cells.Item(141, 9).Value = 9633.333000000001  # CUL!I141
$ws.Cells.Item(141, 11).Value = 28899.999  # CUL!K141
$ws.Cells.Item(141, 13).Value = -23719.999  # CUL!M141

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(123, 8).Value = 33265.832  # GSM!H123
$ws.Cells.Item(123, 10).Value = 33265.832  # GSM!J123
$ws.Cells.Item(123, 12).Value = 33265.832  # GSM!L123
$ws.Cells.Item(123, 14).Value = -38165.832  # GSM!N123

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 8951.888999999999  # LTW!H7
$ws.Cells.Item(7, 9).Value = 9295.357  # LTW!I7
$ws.Cells.Item(7, 11).Value = 9295.357  # LTW!K7
$ws.Cells.Item(7, 13).Value = -9183.357  # LTW!M7

$ws.Cells.Item(61, 8).Value = 3406.8333  # LTW!H61
$ws.Cells.Item(61, 9).Value = 3519  # LTW!I61
$ws.Cells.Item(61, 11).Value = 3519  # LTW!K61
$ws.Cells.Item(61, 13).Value = -3317  # LTW!M61

$ws.Cells.Item(113, 8).Value = 3406.8333  # LTW!H113
$ws.Cells.Item(113, 9).Value = 3519  # LTW!I113
$ws.Cells.Item(113, 11).Value = 3519  # LTW!K113
$ws.Cells.Item(113, 13).Value = -1349  # LTW!M113

$ws.Cells.Item(126, 8).Value = 8951.888999999999  # LTW!H126
$ws.Cells.Item(126, 9).Value = 9295.357  # LTW!I126
$ws.Cells.Item(126, 11).Value = 27886.071  # LTW!K126
$ws.Cells.Item(126, 13).Value = -25416.071  # LTW!M126

$ws.Cells.Item(132, 8).Value = 3279.2942  # LTW!H132
$ws.Cells.Item(132, 9).Value = 2810.8965  # LTW!I132
$ws.Cells.Item(132, 11).Value = 8432.6895  # LTW!K132
$ws.Cells.Item(132, 13).Value = -5902.6895  # LTW!M132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(74, 8).Value = 13359.2  # WVR!H74
$ws.Cells.Item(74, 10).Value = 13359.2  # WVR!J74
$ws.Cells.Item(74, 12).Value = 13359.2  # WVR!L74
$ws.Cells.Item(74, 14).Value = -15231.2  # WVR!N74

$ws.Cells.Item(77, 8).Value = 13359.2  # WVR!H77
$ws.Cells.Item(77, 10).Value = 13359.2  # WVR!J77
$ws.Cells.Item(77, 12).Value = 40077.60000000001  # WVR!L77
$ws.Cells.Item(77, 14).Value = -49437.60000000001  # WVR!N77

$ws.Cells.Item(96, 8).Value = 1232.5555  # WVR!H96
$ws.Cells.Item(96, 9).Value = 973.5  # WVR!I96
$ws.Cells.Item(96, 11).Value = 973.5  # WVR!K96
$ws.Cells.Item(96, 13).Value = 399.5  # WVR!M96

$ws.Cells.Item(105, 8).Value = 0  # WVR!H105
$ws.Cells.Item(105, 10).Value = 0  # WVR!J105
$ws.Cells.Item(105, 12).Value = 0  # WVR!L105
$ws.Cells.Item(105, 14).ClearContents()  # WVR!N105

$ws.Cells.Item(122, 8).Value = 2698.611  # WVR!H122
$ws.Cells.Item(122, 9).Value = 2672.3914  # WVR!I122
$ws.Cells.Item(122, 11).Value = 8017.174199999999  # WVR!K122
$ws.Cells.Item(122, 13).Value = -5567.174199999999  # WVR!M122

$ws.Cells.Item(132, 8).Value = 3528.84  # WVR!H132
$ws.Cells.Item(132, 9).Value = 3341.7046  # WVR!I132
$ws.Cells.Item(132, 11).Value = 10025.1138  # WVR!K132
$ws.Cells.Item(132, 13).Value = -7495.113799999999  # WVR!M132
